$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.387.77"
$ws.Range("E2").Value = "  +7.75%  "
$ws.Range("D3").Value = "3.661.88"
$ws.Range("E3").Value = "  +19.30%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'601.78"
$ws.Range("E5").Value = "  +4.75%  "
$ws.Range("D6").Value = "'186.26"
$ws.Range("E6").Value = "  +9.76%  "
$ws.Range("D7").Value = "3.660.26"
$ws.Range("E7").Value = "  +19.39%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +5.36%  "
$ws.Range("D10").Value = "'0.165"
$ws.Range("E10").Value = "  +10.73%  "
$ws.Range("D11").Value = "'6.57"
$ws.Range("E11").Value = "  +4.94%  "
$ws.Range("E12").Value = "  +7.50%  "
$ws.Range("D13").Value = "'40.18"
$ws.Range("E13").Value = "  +12.97%  "
$ws.Range("D14").Value = "'0.0000258"
$ws.Range("E14").Value = "  +8.30%  "
$ws.Range("D15").Value = "4.278.39"
$ws.Range("E15").Value = "  +19.45%  "
$ws.Range("D16").Value = "71.395.05"
$ws.Range("E16").Value = "  +7.82%  "
$ws.Range("D17").Value = "3.672.35"
$ws.Range("E17").Value = "  +19.49%  "
$ws.Range("E18").Value = "  +2.67%  "
$ws.Range("D19").Value = "'7.57"
$ws.Range("E19").Value = "  +9.20%  "
$ws.Range("E20").Value = "  +4.53%  "
$ws.Range("D21").Value = "'514.55"
$ws.Range("E21").Value = "  +6.25%  "
$ws.Range("D22").Value = "'9.33"
$ws.Range("E22").Value = "  +22.09%  "
$ws.Range("D23").Value = "'0.751"
$ws.Range("E23").Value = "  +9.98%  "
$ws.Range("D24").Value = "'88.65"
$ws.Range("E24").Value = "  +7.90%  "
$ws.Range("D25").Value = "'13.58"
$ws.Range("E25").Value = "  +7.69%  "
$ws.Range("E26").Value = "  +10.18%  "
$ws.Range("E27").Value = "  +8.33%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  +14.15%  "
$ws.Range("D30").Value = "'8.27"
$ws.Range("E30").Value = "  +6.24%  "
$ws.Range("D31").Value = "'32.27"
$ws.Range("E31").Value = "  +17.03%  "
$ws.Range("D32").Value = "'0.0000112"
$ws.Range("E32").Value = "  +22.38%  "
$ws.Range("D33").Value = "'2.77"
$ws.Range("E33").Value = "  +7.21%  "
$ws.Range("E34").Value = "  +5.98%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'6.17"
$ws.Range("E36").Value = "  +11.15%  "
$ws.Range("E37").Value = "  +8.88%  "
$ws.Range("E38").Value = "  +12.82%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'2.13"
$ws.Range("E39").Value = "  +9.18%  "
$ws.Range("B40").Value = "Arweave"
$ws.Range("C40").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D40").Value = "'47.36"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("D41").Value = "'50.96"
$ws.Range("E41").Value = "  +3.94%  "
$ws.Range("E42").Value = "  +5.86%  "
$ws.Range("E43").Value = "  +8.70%  "
$ws.Range("D44").Value = "3.158.66"
$ws.Range("E44").Value = "  +13.89%  "
$ws.Range("D45").Value = "'2.83"
$ws.Range("E45").Value = "  +11.68%  "
$ws.Range("D46").Value = "'407.01"
$ws.Range("E46").Value = "  +11.61%  "
$ws.Range("E47").Value = "  +7.17%  "
$ws.Range("D48").Value = "'28.18"
$ws.Range("E48").Value = "  +16.89%  "
$ws.Range("E49").Value = "  +16.64%  "
$ws.Range("D50").Value = "'134.87"
$ws.Range("E50").Value = "  +0.33%  "
